$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "223.27", "1.00")
# are preserved exactly as text, matching the original inline-string cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '95.767.65'
$ws.Range("E2").Value = '  -1.47%  '

$ws.Range("D3").Value = '3.609.32'
$ws.Range("E3").Value = '  -2.64%  '

$ws.Range("E4").Value = '  +29.71%  '

$ws.Range("E5").Value = '  +0.07%  '

$ws.Range("D6").Value = '223.27'
$ws.Range("E6").Value = '  -5.50%  '

$ws.Range("D7").Value = '638.96'

$ws.Range("D8").Value = '0.424'
$ws.Range("E8").Value = '  -1.75%  '

$ws.Range("D9").Value = '1.22'
$ws.Range("E9").Value = '  +11.27%  '

$ws.Range("E10").Value = '  +0.08%  '

$ws.Range("D11").Value = '3.606.02'
$ws.Range("E11").Value = '  -2.72%  '

$ws.Range("D12").Value = '48.42'
$ws.Range("E12").Value = '  +7.00%  '

$ws.Range("E13").Value = '  +3.25%  '

$ws.Range("D14").Value = '0.0000291'
$ws.Range("E14").Value = '  -5.58%  '

$ws.Range("D15").Value = '6.52'
$ws.Range("E15").Value = '  -5.01%  '

$ws.Range("D16").Value = '4.281.35'
$ws.Range("E16").Value = '  -2.57%  '

$ws.Range("D17").Value = '95.369.64'
$ws.Range("E17").Value = '  -1.64%  '

$ws.Range("D18").Value = '23.83'
$ws.Range("E18").Value = '  +27.33%  '

$ws.Range("D19").Value = '8.95'
$ws.Range("E19").Value = '  -1.99%  '

$ws.Range("D20").Value = '13.79'
$ws.Range("E20").Value = '  +5.78%  '

$ws.Range("D21").Value = '3.606.64'
$ws.Range("E21").Value = '  -2.81%  '

$ws.Range("D22").Value = '0.289'
$ws.Range("E22").Value = '  +41.63%  '

$ws.Range("E23").Value = '  +1.84%  '

$ws.Range("D24").Value = '519.06'
$ws.Range("E24").Value = '  -1.09%  '

$ws.Range("D25").Value = '131.47'
$ws.Range("E25").Value = '  +20.88%  '

$ws.Range("D26").Value = '3.25'
$ws.Range("E26").Value = '  -6.25%  '

$ws.Range("E27").Value = '  -9.19%  '

$ws.Range("D28").Value = '6.79'
$ws.Range("E28").Value = '  -1.40%  '

$ws.Range("D29").Value = '3.775.90'
$ws.Range("E29").Value = '  -3.25%  '

$ws.Range("D30").Value = '12.82'
$ws.Range("E30").Value = '  -5.01%  '

$ws.Range("D31").Value = '13.19'
$ws.Range("E31").Value = '  +4.19%  '

$ws.Range("E32").Value = '  +1.55%  '

$ws.Range("E33").Value = '  +0.12%  '

$ws.Range("D34").Value = '0.627'
$ws.Range("E34").Value = '  +4.99%  '

$ws.Range("D35").Value = '0.182'
$ws.Range("E35").Value = '  -3.84%  '

$ws.Range("D36").Value = '32.95'
$ws.Range("E36").Value = '  +0.75%  '

$ws.Range("E37").Value = '  +0.32%  '

$ws.Range("E38").Value = '  -2.94%  '

$ws.Range("B39").Value = 'USDe'
$ws.Range("C39").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  -0.02%  '

$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '0.535'
$ws.Range("E40").Value = '  +7.55%  '

$ws.Range("D41").Value = '7.28'
$ws.Range("E41").Value = '  +6.69%  '

$ws.Range("D42").Value = '586.13'
$ws.Range("E42").Value = '  -8.38%  '

$ws.Range("D43").Value = '8.34'
$ws.Range("E43").Value = '  -4.62%  '

$ws.Range("E44").Value = '  +14.55%  '

$ws.Range("D45").Value = '41.46'
$ws.Range("E45").Value = '  +3.47%  '

$ws.Range("D46").Value = '0.157'
$ws.Range("E46").Value = '  -6.26%  '

$ws.Range("D47").Value = '0.964'
$ws.Range("E47").Value = '  +0.37%  '

$ws.Range("D48").Value = '1.95'
$ws.Range("E48").Value = '  -4.41%  '

$ws.Range("E49").Value = '  +5.22%  '

$ws.Range("D50").Value = '233.77'
$ws.Range("E50").Value = '  +13.67%  '

$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").Value = '2.27'
$ws.Range("E51").Value = '  -5.15%  '

# Restore the default "Normal" style on column D so no stray number format is left
# behind on cells (matches original styling where these cells had no explicit style).
$ws.Range("D2:D51").Style = "Normal"